# Generate Report for Handoff
# A new handoff round was generated for file "68e85b75-1712-45ed-baff-f14cc6787fa2.md".
# This updates the "Latest Handoff Datetime" for that file on both locale sheets,
# and the "Latest HO Xliff Generate Date" summary on the Overview sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 corresponds to 68e85b75-1712-45ed-baff-f14cc6787fa2.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-26 08:45:54"

# zh-cn sheet: row 7 is the same file, "Latest Handoff Datetime" column H
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-26 08:45:50"

# de-de sheet: row 7 is the same file, "Latest Handoff Datetime" column H
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-26 08:45:54"
